$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ランサーズ")

# Existing hyperlinks do not shift automatically when rows are inserted, so
# clear them first and rebuild the full F2:F26 set afterwards.
$ws1.Hyperlinks.Delete()

# Insert 5 new rows right after the header, pushing the existing listings down.
$ws1.Range("A2:A6").EntireRow.Insert()

# Populate the 5 freshly-scraped listings (newest snapshot, 2025-08-28 12:35:27).
# row 2
$ws1.Range("A2").Value = '2025-08-28 12:35:27'
$ws1.Range("B2").Value = '時給2000円|make/n8n/AWSによる自動化開発者を募集(RPA・MCP経験者歓迎)'
$ws1.Range("C2").Value = 'システム開発'
$ws1.Range("D2").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws1.Range("E2").Value = '期限情報なし'
$ws1.Range("F2").Value = 'https://www.lancers.jp/work/detail/5381595'
$ws1.Range("G2").Value = 143
$ws1.Range("H2").Value = '◆開発,自動化'

# row 3
$ws1.Range("A3").Value = '2025-08-28 12:35:27'
$ws1.Range("B3").Value = '【急募】愛知県弥富市でAMRとPLCを繋ぐシステム開発'
$ws1.Range("C3").Value = 'システム開発'
$ws1.Range("D3").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws1.Range("E3").Value = '期限情報なし'
$ws1.Range("F3").Value = 'https://www.lancers.jp/work/detail/5381608'
$ws1.Range("G3").Value = 125
$ws1.Range("H3").Value = '◆開発,システム開発'

# row 4
$ws1.Range("A4").Value = '2025-08-28 12:35:27'
$ws1.Range("B4").Value = '【急募】ジャストDB(ノーコード)受注システムの編集依頼'
$ws1.Range("C4").Value = 'システム開発'
$ws1.Range("D4").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws1.Range("E4").Value = '期限情報なし'
$ws1.Range("F4").Value = 'https://www.lancers.jp/work/detail/5379679'
$ws1.Range("G4").Value = 33

# row 5
$ws1.Range("A5").Value = '2025-08-28 12:35:27'
$ws1.Range("B5").Value = 'amazon注文詳細より配送番号と配送業者と配送状況取得してもらいたいです。'
$ws1.Range("C5").Value = 'システム開発'
$ws1.Range("D5").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws1.Range("E5").Value = '期限情報なし'
$ws1.Range("F5").Value = 'https://www.lancers.jp/work/detail/5381625'
$ws1.Range("G5").Value = 25

# row 6
$ws1.Range("A6").Value = '2025-08-28 12:35:27'
$ws1.Range("B6").Value = '【急募】利用者予定表と勤務表の自動集計をお任せ!'
$ws1.Range("C6").Value = 'システム開発'
$ws1.Range("D6").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws1.Range("E6").Value = '期限情報なし'
$ws1.Range("F6").Value = 'https://www.lancers.jp/work/detail/5381634'
$ws1.Range("G6").Value = 13

# Re-create the hyperlinks on column F for every data row (2-26), pointing at
# the listing URL already written into that cell.
for ($r = 2; $r -le 26; $r++) {
    $cell = $ws1.Range("F" + $r)
    $ws1.Hyperlinks.Add($cell, $cell.Value2)
}

# --- 統計 (stats) sheet: append the summary row for this scrape run ---
$ws2 = $wb.Worksheets.Item("統計")
$ws2.Range("A11").Value = '2025-08-28T12:35:27.199383'
$ws2.Range("B11").Value = 13
$ws2.Range("C11").Value = '全案件リスト'
$ws2.Range("D11").Value = 53.8
$ws2.Range("E11").Value = 6
$ws2.Range("F11").Value = 4
$ws2.Range("G11").Value = 13
